$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44511
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25375
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("S2").Value = 2538
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44483
$ws.Range("L3").Value = 'Primera'

# Row 4
$ws.Range("D4").Value = 44483
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 2400
$ws.Range("O4").Value = 2400
$ws.Range("P4").Value = 2400
$ws.Range("Q4").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S4").Value = 2400
$ws.Range("T4").Value = 1

# Row 5
$ws.Range("D5").Value = 44467
$ws.Range("N5").Value = 2700
$ws.Range("O5").Value = 2800
$ws.Range("P5").Value = 2750
$ws.Range("Q5").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S5").Value = 2750
$ws.Range("T5").Value = 1

# Row 6
$ws.Range("D6").Value = 44467
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 2500
$ws.Range("O6").Value = 2500
$ws.Range("P6").Value = 2500
$ws.Range("Q6").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S6").Value = 2500
$ws.Range("T6").Value = 1

# Row 8
$ws.Range("D8").Value = 44846
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 23000
$ws.Range("O8").Value = 23000
$ws.Range("P8").Value = 23000
$ws.Range("Q8").Value = '$/bandeja 10 kilos'
$ws.Range("S8").Value = 2300
$ws.Range("T8").Value = 10

# Row 9
$ws.Range("D9").Value = 44846
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("Q9").Value = '$/bandeja 10 kilos'
$ws.Range("S9").Value = 2050
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44804
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 29000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 29500
$ws.Range("S10").Value = 2950

# Row 11
$ws.Range("D11").Value = 44446
$ws.Range("N11").Value = 3200
$ws.Range("O11").Value = 3300
$ws.Range("P11").Value = 3250
$ws.Range("R11").Value = 'Provincia del Elquí'
$ws.Range("S11").Value = 3250

# Row 12
$ws.Range("D12").Value = 44461

# Row 13
$ws.Range("D13").Value = 44530
$ws.Range("N13").Value = 2000
$ws.Range("O13").Value = 2100
$ws.Range("P13").Value = 2050
$ws.Range("S13").Value = 2050

# Row 14
$ws.Range("D14").Value = 44819
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 80
$ws.Range("N14").Value = 25000
$ws.Range("O14").Value = 26000
$ws.Range("P14").Value = 25500
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("S14").Value = 2550
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44160
$ws.Range("L15").Value = 'Primera'
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 18000
$ws.Range("P15").Value = 17500
$ws.Range("Q15").Value = '$/bandeja 8 kilos'
$ws.Range("S15").Value = 2188
$ws.Range("T15").Value = 8

# Row 16
$ws.Range("D16").Value = 44160
$ws.Range("L16").Value = 'Segunda'
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("Q16").Value = '$/bandeja 8 kilos'
$ws.Range("S16").Value = 1875
$ws.Range("T16").Value = 8

# Row 17
$ws.Range("D17").Value = 44491
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = 25000
$ws.Range("O17").Value = 26000
$ws.Range("P17").Value = 25467
$ws.Range("S17").Value = 2547

# Row 18
$ws.Range("D18").Value = 44516
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 1900
$ws.Range("O18").Value = 2000
$ws.Range("P18").Value = 1950
$ws.Range("Q18").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S18").Value = 1950
$ws.Range("T18").Value = 1

# Row 19
$ws.Range("D19").Value = 44516
$ws.Range("L19").Value = 'Segunda'
$ws.Range("N19").Value = 1700
$ws.Range("O19").Value = 1700
$ws.Range("P19").Value = 1700
$ws.Range("S19").Value = 1700

# Row 20
$ws.Range("D20").Value = 44469
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 28000
$ws.Range("O20").Value = 29000
$ws.Range("P20").Value = 28500
$ws.Range("Q20").Value = '$/bandeja 10 kilos'
$ws.Range("S20").Value = 2850
$ws.Range("T20").Value = 10

# Row 21
$ws.Range("D21").Value = 44505
$ws.Range("N21").Value = 2200
$ws.Range("O21").Value = 2200
$ws.Range("P21").Value = 2200
$ws.Range("Q21").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S21").Value = 2200
$ws.Range("T21").Value = 1

# Row 22
$ws.Range("D22").Value = 44505
$ws.Range("L22").Value = 'Segunda'
$ws.Range("N22").Value = 1800
$ws.Range("O22").Value = 1800
$ws.Range("P22").Value = 1800
$ws.Range("S22").Value = 1800

# Row 23
$ws.Range("D23").Value = 44806
$ws.Range("L23").Value = 'Especial'
$ws.Range("N23").Value = 22000
$ws.Range("O23").Value = 22000
$ws.Range("P23").Value = 22000
$ws.Range("Q23").Value = '$/bandeja 10 kilos'
$ws.Range("S23").Value = 2200
$ws.Range("T23").Value = 10

# Row 24
$ws.Range("D24").Value = 44806
$ws.Range("L24").Value = 'Primera'
$ws.Range("N24").Value = 19000
$ws.Range("O24").Value = 19000
$ws.Range("P24").Value = 19000
$ws.Range("S24").Value = 1900

# Row 25
$ws.Range("D25").Value = 44806
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("S25").Value = 1700

# Row 26
$ws.Range("D26").Value = 44454
$ws.Range("N26").Value = 30000
$ws.Range("O26").Value = 31000
$ws.Range("P26").Value = 30500
$ws.Range("Q26").Value = '$/bandeja 10 kilos'
$ws.Range("S26").Value = 3050
$ws.Range("T26").Value = 10

# Row 27
$ws.Range("D27").Value = 44488
$ws.Range("L27").Value = 'Primera'
$ws.Range("N27").Value = 25000
$ws.Range("O27").Value = 26000
$ws.Range("P27").Value = 25600
$ws.Range("Q27").Value = '$/bandeja 10 kilos'
$ws.Range("S27").Value = 2560
$ws.Range("T27").Value = 10

# Row 28
$ws.Range("D28").Value = 44517
$ws.Range("L28").Value = 'Primera'
$ws.Range("N28").Value = 25000
$ws.Range("O28").Value = 27000
$ws.Range("P28").Value = 26000
$ws.Range("S28").Value = 2600

# Row 29
$ws.Range("D29").Value = 44462
$ws.Range("N29").Value = 2900
$ws.Range("O29").Value = 3000
$ws.Range("P29").Value = 2950
$ws.Range("Q29").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S29").Value = 2950
$ws.Range("T29").Value = 1

# Row 30
$ws.Range("D30").Value = 44462
$ws.Range("L30").Value = 'Segunda'
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 2600
$ws.Range("O30").Value = 2600
$ws.Range("P30").Value = 2600
$ws.Range("Q30").Value = '$/kilo (en caja de 15 kilos)'
$ws.Range("S30").Value = 2600
$ws.Range("T30").Value = 1

# Row 32
$ws.Range("D32").Value = 44495
$ws.Range("N32").Value = 26000
$ws.Range("O32").Value = 27000
$ws.Range("P32").Value = 26500
$ws.Range("Q32").Value = '$/bandeja 10 kilos'
$ws.Range("R32").Value = 'Provincia de Limarí'
$ws.Range("S32").Value = 2650
$ws.Range("T32").Value = 10

# Row 33
$ws.Range("D33").Value = 44832
$ws.Range("N33").Value = 25000
$ws.Range("O33").Value = 26000
$ws.Range("P33").Value = 25500
$ws.Range("S33").Value = 2550

# Row 34
$ws.Range("D34").Value = 44845
$ws.Range("L34").Value = 'Especial'
$ws.Range("N34").Value = 23000
$ws.Range("O34").Value = 23000
$ws.Range("P34").Value = 23000
$ws.Range("S34").Value = 2300

# Row 35
$ws.Range("D35").Value = 44845
$ws.Range("N35").Value = 21000
$ws.Range("O35").Value = 21000
$ws.Range("P35").Value = 21000
$ws.Range("S35").Value = 2100
